# Mise a jour mapping posologie pour reunion ef2bdf36d95fd7302ceb15e853eb11298b18e179
#
# 1) Bump the Metadata "Date" value.
# 2) Extend "Mapping Table 2" with 5 more rows that repeat the existing
#    Frq_filtreVal_1_J mapping row (rows 5-9).
# 3) Add a new "Mapping Table 3" sheet (same layout as the other mapping
#    tables) describing the Quantite/Nombre and Quantite/Unite mappings.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 date bump -------------------------------------------
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2024-11-18T17:53:13+00:00"

# --- 2) Grow "Mapping Table 2" from 4 to 9 rows --------------------------
$table2 = $wb.Worksheets.Item("Mapping Table 2")
$patternRow = $table2.Range("A3:E3")

for ($r = 5; $r -le 9; $r++) {
    $destRow = $table2.Range("A" + $r + ":E" + $r)
    $patternRow.Copy()
    $destRow.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $patternRow.Copy()
    $destRow.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# --- 3) New "Mapping Table 3" sheet --------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$table3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$table3.Name = "Mapping Table 3"

# Header row + blank spacer row, copied (values + formats) from Table 2.
$headerSrc = $table2.Range("A1:E2")
$headerDest = $table3.Range("A1:E2")
$headerSrc.Copy()
$headerDest.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$headerSrc.Copy()
$headerDest.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row formatting for the two data rows, copied from Table 2 rows 3-4.
$dataFmtSrc = $table2.Range("A3:E4")
$dataFmtDest = $table3.Range("A3:E4")
$dataFmtSrc.Copy()
$dataFmtDest.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$table3.Range("A3").Value = "Elément_posologie/Quantité/Nombre"
$table3.Range("C3").Value = "equivalent"
$table3.Range("D3").Value = "Dosage.doseAndRate.doseQuantity.value"

$table3.Range("A4").Value = "Elément_posologie/Quantité/Unité"
$table3.Range("C4").Value = "related-to"
$table3.Range("D4").Value = "Dosage.doseAndRate.doseQuantity.unité"

# Restore original active sheet/tab selection (Metadata, tab 0).
$metaSheet.Activate()
